# Generate Report for Handback
#
# This mirrors the "handback" report-generation step: the localized
# files (zh-cn, de-de) have been handed back, so for both rows on each
# language sheet we now know the Latest Target File (the source .md,
# now hyperlinked the same way the Source File Name column is), the
# Latest Handback File (the .xlf that was handed back) and the real
# Latest Handback DateTime (replacing the zero-date placeholder). The
# Status column (and the Overview roll-up columns) flips from
# "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$ov   = $wb.Worksheets.Item("Overview")
$zh   = $wb.Worksheets.Item("zh-cn")
$de   = $wb.Worksheets.Item("de-de")

$status = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de roll-up status columns ---
$ov.Range("E2").Value = $status
$ov.Range("F2").Value = $status
$ov.Range("E3").Value = $status
$ov.Range("F3").Value = $status

# --- zh-cn sheet ---
$zh.Range("C2").Value = $status
$zh.Range("C3").Value = $status

$zh.Range("I2").Value = "416f48b4-653f-46a1-9368-33089e2614cc.md"
$zh.Range("J2").Value = "416f48b4-653f-46a1-9368-33089e2614cc.624efe40e85c3a0a91593478ec01b79f3ac83a7c.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-03 04:52:37"

$zh.Range("I3").Value = "a6869899-3547-46eb-ba75-ae8e59a522d9.md"
$zh.Range("J3").Value = "a6869899-3547-46eb-ba75-ae8e59a522d9.03e74bf2ceace2fa73293d73198009e554b79884.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-03 04:52:37"

$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b23aea7a43fc87d83ddeaba54a15383b0763d8e9/e2e/416f48b4-653f-46a1-9368-33089e2614cc.md", "", "", "416f48b4-653f-46a1-9368-33089e2614cc.md")
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b23aea7a43fc87d83ddeaba54a15383b0763d8e9/e2e/a6869899-3547-46eb-ba75-ae8e59a522d9.md", "", "", "a6869899-3547-46eb-ba75-ae8e59a522d9.md")

# Match the workbook's existing custom HyperLink look (cornflower blue,
# underlined) instead of Excel's default theme hyperlink colour.
$zh.Range("I2").Font.Underline = $true
$zh.Range("I2").Font.Color = 15570276
$zh.Range("I3").Font.Underline = $true
$zh.Range("I3").Font.Color = 15570276

# --- de-de sheet ---
$de.Range("C2").Value = $status
$de.Range("C3").Value = $status

$de.Range("I2").Value = "416f48b4-653f-46a1-9368-33089e2614cc.md"
$de.Range("J2").Value = "416f48b4-653f-46a1-9368-33089e2614cc.624efe40e85c3a0a91593478ec01b79f3ac83a7c.de-de.xlf"
$de.Range("K2").Value = "2016-09-03 04:52:44"

$de.Range("I3").Value = "a6869899-3547-46eb-ba75-ae8e59a522d9.md"
$de.Range("J3").Value = "a6869899-3547-46eb-ba75-ae8e59a522d9.03e74bf2ceace2fa73293d73198009e554b79884.de-de.xlf"
$de.Range("K3").Value = "2016-09-03 04:52:44"

$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b23aea7a43fc87d83ddeaba54a15383b0763d8e9/e2e/416f48b4-653f-46a1-9368-33089e2614cc.md", "", "", "416f48b4-653f-46a1-9368-33089e2614cc.md")
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b23aea7a43fc87d83ddeaba54a15383b0763d8e9/e2e/a6869899-3547-46eb-ba75-ae8e59a522d9.md", "", "", "a6869899-3547-46eb-ba75-ae8e59a522d9.md")

# Match the workbook's existing custom HyperLink look (cornflower blue,
# underlined) instead of Excel's default theme hyperlink colour.
$de.Range("I2").Font.Underline = $true
$de.Range("I2").Font.Color = 15570276
$de.Range("I3").Font.Underline = $true
$de.Range("I3").Font.Color = 15570276

# --- Column widths: columns widened to fit the now-populated
#     Status / Latest Target File / Latest Handback File columns ---
$ov.Columns.Item(5).AutoFit()
$ov.Columns.Item(6).AutoFit()

$zh.Columns.Item(3).AutoFit()
$zh.Columns.Item(9).AutoFit()
$zh.Columns.Item(10).AutoFit()

$de.Columns.Item(3).AutoFit()
$de.Columns.Item(9).AutoFit()
$de.Columns.Item(10).AutoFit()
